# Clean up stray test/debug data left in the "Quater" columns and fix a
# couple of student-name typos, then leave the selection where the user
# was last working (so the app can show them where the fix landed).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Fix mistyped / duplicated student names -------------------------------
$ws.Range("B6").Value = "Arissara Wuthisasna"
$ws.Range("B8").Value = "Chakree Kenganantanon"

# --- Replace the leftover scratch/test values in columns C:F ---------------
$ws.Range("D5").ClearContents()
$ws.Range("C5").Value = "test"

$ws.Range("D6").Value = "test"

$ws.Range("C7").ClearContents()
$ws.Range("D7").ClearContents()
$ws.Range("E7").ClearContents()

$ws.Range("D8").Value = "tse"
$ws.Range("E8").ClearContents()
$ws.Range("F8").ClearContents()

$ws.Range("D10").ClearContents()
$ws.Range("D11").ClearContents()

$ws.Range("E13").ClearContents()

# --- Leave the selection on the cell the user last touched -----------------
$ws.Range("C19").Select()
